$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.286.05"
$ws.Range("E2").Value = "  +0.70%  "

$ws.Range("D3").Value = "1.664.64"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("D5").Value = "218.76"
$ws.Range("E5").Value = "  +0.47%  "

$ws.Range("D6").Value = "0.5320"
$ws.Range("E6").Value = "  +1.33%  "

$ws.Range("E7").Value = "  +0.81%  "

$ws.Range("D8").Value = "0.2642"
$ws.Range("E8").Value = "  +1.32%  "

$ws.Range("E9").Value = "  +0.45%  "

$ws.Range("E10").Value = "  +0.78%  "

$ws.Range("D11").Value = "0.07825"
$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("D12").Value = "4.552"
$ws.Range("E12").Value = "  +0.94%  "

$ws.Range("D13").Value = "1.661.93"
$ws.Range("E13").Value = "  +1.69%  "

$ws.Range("D14").Value = "1.893.10"
$ws.Range("E14").Value = "  +0.52%  "

$ws.Range("D15").Value = "0.5518"
$ws.Range("E15").Value = "  +0.77%  "

$ws.Range("D16").Value = "0.0₅8210"
$ws.Range("E16").Value = "  -0.05%  "

$ws.Range("D17").Value = "65.63"
$ws.Range("E17").Value = "  +0.35%  "

$ws.Range("E18").Value = "  +0.86%  "

$ws.Range("D19").Value = "4.687"
$ws.Range("E19").Value = "  +2.20%  "

$ws.Range("D20").Value = "193.48"
$ws.Range("E20").Value = "  +1.15%  "

$ws.Range("E21").Value = "  +1.35%  "

$ws.Range("D22").Value = "6.028"
$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D24").Value = "145.58"
$ws.Range("E24").Value = "  +2.45%  "

$ws.Range("D25").Value = "0.1228"
$ws.Range("E25").Value = "  -0.93%  "

$ws.Range("D26").Value = "7.196"
$ws.Range("E26").Value = "  -0.80%  "

$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("E28").Value = "  +3.55%  "

$ws.Range("D29").Value = "0.05860"
$ws.Range("E29").Value = "  -0.66%  "

$ws.Range("D30").Value = "1.282"
$ws.Range("E30").Value = "  +0.54%  "

$ws.Range("D31").Value = "3.607"
$ws.Range("E31").Value = "  +2.63%  "

$ws.Range("D32").Value = "3.277"
$ws.Range("E32").Value = "  +0.62%  "

$ws.Range("D33").Value = "1.608"
$ws.Range("E33").Value = "  +1.11%  "

$ws.Range("D34").Value = "0.9621"
$ws.Range("E34").Value = "  +1.15%  "

$ws.Range("D35").Value = "2.827"
$ws.Range("E35").Value = "  +1.62%  "

$ws.Range("D36").Value = "2.418"
$ws.Range("E36").Value = "  +0.38%  "

$ws.Range("D37").Value = "0.5803"
$ws.Range("E37").Value = "  +2.10%  "

$ws.Range("D38").Value = "0.01607"
$ws.Range("E38").Value = "  -0.65%  "

$ws.Range("D39").Value = "0.8661"
$ws.Range("E39").Value = "  +2.07%  "

$ws.Range("D40").Value = "5.829"
$ws.Range("E40").Value = "  +0.32%  "

$ws.Range("D41").Value = "1.049.96"
$ws.Range("E41").Value = "  +1.84%  "

$ws.Range("D42").Value = "1.010"
$ws.Range("E42").Value = "  +0.78%  "

$ws.Range("D43").Value = "104.46"
$ws.Range("E43").Value = "  +1.60%  "

$ws.Range("D44").Value = "1.803.77"
$ws.Range("E44").Value = "  +0.34%  "

$ws.Range("D45").Value = "57.69"
$ws.Range("E45").Value = "  +0.89%  "

$ws.Range("D46").Value = "1.005"
$ws.Range("E46").Value = "  +0.33%  "

$ws.Range("E47").Value = "  -7.60%  "

$ws.Range("D48").Value = "0.4383"
$ws.Range("E48").Value = "  +1.86%  "

$ws.Range("D49").Value = "8.061"
$ws.Range("E49").Value = "  +2.64%  "

$ws.Range("D50").Value = "0.05163"
$ws.Range("E50").Value = "  +0.11%  "

$ws.Range("D51").Value = "1.420"
$ws.Range("E51").Value = "  -3.89%  "
